# Updated cryptos list (Price/Volume(1h) refresh + a rank swap between
# Mantle and InjectiveProtocol in rows 48-49), mirroring the GitHub
# Actions scheduled-update commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matches source inline-string cells) for D-column
# price values that would otherwise be auto-parsed as numbers by Excel.
foreach ($addr in @("D5", "D6", "D8", "D14", "D20", "D21", "D22", "D24", "D25", "D26", "D37", "D38", "D40", "D45", "D46", "D48", "D49")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '62.824.96'
$ws.Range("E2").Value = '  -0.61%  '
$ws.Range("D3").Value = '2.464.44'
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '571.93'
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("D6").Value = '147.45'
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.531'
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("E10").Value = '  +0.03%  '
$ws.Range("E11").Value = '  -1.25%  '
$ws.Range("E12").Value = '  -1.54%  '
$ws.Range("E13").Value = '  +1.92%  '
$ws.Range("D14").Value = '0.0000176'
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("D15").Value = '2.908.79'
$ws.Range("E15").Value = '  -0.65%  '
$ws.Range("D16").Value = '62.746.96'
$ws.Range("E16").Value = '  -0.65%  '
$ws.Range("D17").Value = '2.465.64'
$ws.Range("E17").Value = '  -0.23%  '
$ws.Range("E18").Value = '  -5.84%  '
$ws.Range("E19").Value = '  -2.57%  '
$ws.Range("D20").Value = '2.34'
$ws.Range("E20").Value = '  +4.38%  '
$ws.Range("D21").Value = '4.15'
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("D22").Value = '321.41'
$ws.Range("E22").Value = '  -2.58%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = '10.18'
$ws.Range("E24").Value = '  +2.57%  '
$ws.Range("D25").Value = '64.77'
$ws.Range("E25").Value = '  -2.16%  '
$ws.Range("D26").Value = '640.21'
$ws.Range("E26").Value = '  -2.34%  '
$ws.Range("D28").Value = '0.0₃0965'
$ws.Range("E28").Value = '  -2.64%  '
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("E30").Value = '  -4.12%  '
$ws.Range("E31").Value = '  -2.01%  '
$ws.Range("E32").Value = '  -2.05%  '
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("E35").Value = '  -3.19%  '
$ws.Range("E36").Value = '  -2.16%  '
$ws.Range("D37").Value = '5.37'
$ws.Range("E37").Value = '  -1.51%  '
$ws.Range("D38").Value = '0.366'
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("E39").Value = '  -1.26%  '
$ws.Range("D40").Value = '148.35'
$ws.Range("E40").Value = '  -1.32%  '
$ws.Range("E41").Value = '  -1.21%  '
$ws.Range("E42").Value = '  -1.73%  '
$ws.Range("D43").Value = '0.0₆0308'
$ws.Range("E43").Value = '  -3.72%  '
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").Value = '154.35'
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("D46").Value = '15.38'
$ws.Range("E46").Value = '  +0.96%  '
$ws.Range("E47").Value = '  -1.34%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '0.606'
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '20.24'
$ws.Range("E49").Value = '  -0.90%  '
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("E51").Value = '  -1.72%  '
